# Apply the cryptos list update for Sat Jan 20 22:59:16 UTC 2024 (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "41.768.78"
$ws.Range("E2").Value = "  +0.57%  "

# Row 3
$ws.Range("D3").Value = "2.471.76"
$ws.Range("E3").Value = "  -0.70%  "

# Row 4
$ws.Range("E4").Value = "  +0.14%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.14%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "93.15"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.25%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.550"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.60%  "

# Row 8
$ws.Range("E8").Value = "  +0.12%  "

# Row 9
$ws.Range("E9").Value = "  +3.41%  "

# Row 10
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0863"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +10.09%  "

# Row 11
$ws.Range("B11").Value = "Avalanche"
$ws.Range("C11").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "32.93"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.87%  "

# Row 12
$ws.Range("E12").Value = "  +0.14%  "

# Row 13
$ws.Range("D13").Value = "2.851.47"
$ws.Range("E13").Value = "  -0.55%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.91"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.92%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.83"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.09%  "

# Row 16
$ws.Range("D16").Value = "2.487.84"
$ws.Range("E16").Value = "  +0.87%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.784"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.95%  "

# Row 18
$ws.Range("D18").Value = "41.741.72"
$ws.Range("E18").Value = "  +0.43%  "

# Row 19
$ws.Range("D19").Value = "0.0₃0963"
$ws.Range("E19").Value = "  +3.74%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.51"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.52%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.57"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.69%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.22"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.10%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "239.82"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.54%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.73"
$ws.Range("D24").Style = "Normal"

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.92"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.76%  "

# Row 26
$ws.Range("E26").Value = "  -0.04%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.79"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.69%  "

# Row 28
$ws.Range("E28").Value = "  +2.22%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.82"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.57%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.51"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.09%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "156.21"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.69%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.52"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.35%  "

# Row 33
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0767"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.55%  "

# Row 34
$ws.Range("B34").Value = "WEMIXToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.58"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.50%  "

# Row 35
$ws.Range("E35").Value = "  +2.53%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.55"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.92%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.90"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.24%  "

# Row 38
$ws.Range("E38").Value = "  +1.10%  "

# Row 39
$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.82"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.59%  "

# Row 40
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.103"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.30%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.84%  "

# Row 42
$ws.Range("E42").Value = "  +0.03%  "

# Row 43
$ws.Range("D43").Value = "1.971.52"
$ws.Range("E43").Value = "  +0.30%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.08"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.01%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0286"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.65%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.95"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.92%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.08"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.47%  "

# Row 48
$ws.Range("D48").Value = "2.704.21"
$ws.Range("E48").Value = "  -0.83%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "97.51"
$ws.Range("D49").Style = "Normal"

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "67.25"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.65%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "52.93"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.79%  "
